$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7234318852424622
$ws.Range("B1").Value = 0.9965925812721252
$ws.Range("C1").Value = 0.8978086113929749
$ws.Range("D1").Value = 3.25597071647644
$ws.Range("E1").Value = 1.619717001914978
